$d = $word.ActiveDocument

# Namespace-qualified pkg wrapper used for all InsertXML calls below. InsertXML
# REPLACES exactly the contents of the Range it is called on, giving full
# control over run-level formatting (w:sz / w:szCs / w:proofErr, etc.) that
# isn't reachable purely through Font.* properties in this host.
function New-PkgXml([string]$bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $bodyInner +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# -----------------------------------------------------------------------
# 1) Title paragraph: "Week 4" -> "Neural Networks: Representation"
#    Also drop the _Toc471159746 bookmark that wrapped the old title.
# -----------------------------------------------------------------------
$d.Bookmarks.Item("_Toc471159746").Delete()

$p = $d.Paragraphs.Item(1)
$titleRange = $p.Range
$titleRange.MoveEnd(1, -1) | Out-Null

$titleRunsXml =
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="80"/><w:szCs w:val="80"/></w:rPr><w:t>N</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>eural</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="80"/><w:szCs w:val="80"/></w:rPr><w:t xml:space="preserve"> N</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>etworks</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="80"/><w:szCs w:val="80"/></w:rPr><w:t xml:space="preserve"> R</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>epresentation</w:t></w:r>'

$titleRange.InsertXML((New-PkgXml $titleRunsXml))

# -----------------------------------------------------------------------
# 2) "O(n^degree) ~ (n^2)/2" - collapse the spell-check-split runs into one
# -----------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("O(n^degree) ~ (n^2)/2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1 = $d.Range($find1.Start, $find1.End)

$r1.InsertXML((New-PkgXml '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>O(n^degree) ~ (n^2)/2</w:t></w:r>'))

# -----------------------------------------------------------------------
# 3) "... x (sj + 1)" - collapse the spell-check-split runs for "s"/"j"
# -----------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("j + 1 x (sj + 1)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2 = $d.Range($find2.Start, $find2.End)

$seg2Xml =
    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>j + 1</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> x (s</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">j </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:val="en-US"/></w:rPr><w:t>+ 1)</w:t></w:r>'

$r2.InsertXML((New-PkgXml $seg2Xml))

Write-Host ("Title=[" + $d.Paragraphs.Item(1).Range.Text + "]")
